$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column B ("customer_id") - shifts old B..L one column right
#    (old B:address -> C, C:npwp -> D, D:coordinate -> E, E:email -> F,
#     F:phone -> G, G:tax_invoice_number -> H, H:package_id -> I,
#     I:status -> J, J:join_date -> K, K:bill_date -> L, L:inactive_at -> M)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Range("B1").Value = "customer_id"

# Customer id values (entered in this order so the shared-string table comes
# out IDPEL-001, IDPEL-002, IDPEL-000, IDPEL-010 - row 2 / "A. Mulyana" is
# left without a customer_id)
$ws.Range("B4").Value = "IDPEL-001"
$ws.Range("B3").Value = "IDPEL-002"
$ws.Range("B5").Value = "IDPEL-000"
$ws.Range("B6").Value = "IDPEL-010"

# Match column A's width for the new column
$ws.Columns.Item(2).ColumnWidth = 25.26

# ---------------------------------------------------------------------------
# 2) bill_date (now column L after the insert) gets a day-of-month number
#    for each customer, formatted as a plain integer
# ---------------------------------------------------------------------------
$ws.Range("L1:L6").NumberFormat = "0"
$ws.Range("L2").Value = 10
$ws.Range("L3").Value = 22
$ws.Range("L4").Value = 20
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 5

$ws.Columns.Item(12).ColumnWidth = 10.17

# ---------------------------------------------------------------------------
# 3) Update the view: scroll so column F is leftmost, select M7
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("F1"), $true)
$ws.Range("M7").Select()
